$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (CB_API)
$ws.Range("S2").Value = 39
$ws.Range("T2").Value = 183.5
$ws.Range("U2").Value = 512.76
$ws.Range("V2").Value = 980.15

# Row 14 updates (Dash)
$ws.Range("J14").Value = 1123.31
$ws.Range("K14").Value = 1416.81
$ws.Range("L14").Value = 2263.59
$ws.Range("M14").Value = 2047.83
$ws.Range("N14").Value = 886.8100000000001
$ws.Range("O14").Value = 989.91
$ws.Range("P14").Value = 1318.75
$ws.Range("Q14").Value = 1264.24
$ws.Range("R14").Value = 1710.87
$ws.Range("S14").Value = 2117.04
$ws.Range("T14").Value = 2918.95
$ws.Range("U14").Value = 1775.78
$ws.Range("V14").Value = 1273.37
$ws.Range("W14").Value = 428.78
